# Add a new day ("2025-12-29") to the "Chart" sheet, mirroring the existing
# breadcrumb rows (Date string in column A, Invalid count in B, Valid count in C).
#
# Column A stores the date as plain text (shared string), not a real Excel
# date. Assigning a date-like string straight to .Value makes Excel helpfully
# (but unhelpfully for us) re-interpret it as a date serial number. To avoid
# that, we build the text in a scratch cell with a non-date-like trailing
# space (so Excel leaves it alone as text), clean the trailing space with a
# TRIM() formula, and paste the *value* of that formula into place. That
# keeps the destination cell a plain shared-string text cell using the
# workbook's default style, just like all its neighboring date cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$newRow = 86
$dateText = "2025-12-29"

# --- Column A: the date, stored as text -----------------------------------
$scratch1 = $ws.Range("Z1")
$scratch2 = $ws.Range("Z2")

$scratch1.Value = $dateText + " "
$scratch2.Formula = "=TRIM(Z1)"

$scratch2.Copy()
$ws.Range("A" + $newRow).PasteSpecial(-4163)

$scratch1.Clear()
$scratch2.Clear()

# --- Columns B & C: numeric counts ------------------------------------------
$ws.Range("B" + $newRow).Value = 0
$ws.Range("C" + $newRow).Value = 28

Write-Host "Added row $newRow with date $dateText"
